$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 18.24971631898057
$ws.Cells.Item(2, 3).Value = 9.131796674211987
$ws.Cells.Item(2, 4).Value = 5.993911966954698
$ws.Cells.Item(2, 5).Value = 11.16960707678143
$ws.Cells.Item(2, 7).Value = 3.650426866380956
$ws.Cells.Item(2, 12).Value = 9.908986248679096
$ws.Cells.Item(2, 14).Value = 17.92730646606201
$ws.Cells.Item(2, 15).Value = 25.70840185818437
$ws.Cells.Item(3, 2).Value = 17.76209409227869
$ws.Cells.Item(3, 3).Value = 8.86726559832891
$ws.Cells.Item(3, 4).Value = 5.877203648217691
$ws.Cells.Item(3, 5).Value = 11.20378971940537
$ws.Cells.Item(3, 7).Value = 3.653170439028423
$ws.Cells.Item(3, 12).Value = 9.886328083618078
$ws.Cells.Item(3, 14).Value = 17.99481942739907
$ws.Cells.Item(3, 15).Value = 25.71293567331942
$ws.Cells.Item(4, 2).Value = 17.45957413723468
$ws.Cells.Item(4, 3).Value = 8.699274718298364
$ws.Cells.Item(4, 4).Value = 5.806216850924396
$ws.Cells.Item(4, 5).Value = 11.22652872423583
$ws.Cells.Item(4, 7).Value = 3.654944129145282
$ws.Cells.Item(4, 12).Value = 9.874260082675958
$ws.Cells.Item(4, 14).Value = 18.03820946756112
$ws.Cells.Item(4, 15).Value = 25.72300375325649
$ws.Cells.Item(5, 2).Value = 17.3357049778715
$ws.Cells.Item(5, 3).Value = 8.629479341981547
$ws.Cells.Item(5, 4).Value = 5.777502238205807
$ws.Cells.Item(5, 5).Value = 11.23623514668801
$ws.Cells.Item(5, 7).Value = 3.655689410017274
$ws.Cells.Item(5, 12).Value = 9.869809356184591
$ws.Cells.Item(5, 14).Value = 18.05637985636202
$ws.Cells.Item(5, 15).Value = 25.72893462336969
$ws.Cells.Item(6, 2).Value = 17.31510683947879
$ws.Cells.Item(6, 3).Value = 8.617811004393158
$ws.Cells.Item(6, 4).Value = 5.772748357039532
$ws.Cells.Item(6, 5).Value = 11.23787346720531
$ws.Cells.Item(6, 7).Value = 3.655814523790864
$ws.Cells.Item(6, 12).Value = 9.869098616883395
$ws.Cells.Item(6, 14).Value = 18.05942659173692
$ws.Cells.Item(6, 15).Value = 25.73002971880686
$ws.Cells.Item(7, 2).Value = 17.45790571887084
$ws.Cells.Item(7, 3).Value = 8.698338763539345
$ws.Cells.Item(7, 4).Value = 5.805828674587278
$ws.Cells.Item(7, 5).Value = 11.22665784670602
$ws.Cells.Item(7, 7).Value = 3.654954089110935
$ws.Cells.Item(7, 12).Value = 9.874198163306563
$ws.Cells.Item(7, 14).Value = 18.0384525394144
$ws.Cells.Item(7, 15).Value = 25.72307634315452
$ws.Cells.Item(8, 2).Value = 18.08234472201384
$ws.Cells.Item(8, 3).Value = 9.041777665073141
$ws.Cells.Item(8, 4).Value = 5.953557616243973
$ws.Cells.Item(8, 5).Value = 11.18102980599603
$ws.Cells.Item(8, 7).Value = 3.651354398335009
$ws.Cells.Item(8, 12).Value = 9.900793069767479
$ws.Cells.Item(8, 14).Value = 17.95018389723532
$ws.Cells.Item(8, 15).Value = 25.70845142469641
$ws.Cells.Item(9, 2).Value = 19.27412443786747
$ws.Cells.Item(9, 3).Value = 9.668617931840421
$ws.Cells.Item(9, 4).Value = 6.246683048012276
$ws.Cells.Item(9, 5).Value = 11.10545085208855
$ws.Cells.Item(9, 7).Value = 3.644999100672438
$ws.Cells.Item(9, 12).Value = 9.967415219573718
$ws.Cells.Item(9, 14).Value = 17.79238473167417
$ws.Cells.Item(9, 15).Value = 25.7376802844897
$ws.Cells.Item(10, 2).Value = 20.11985487985555
$ws.Cells.Item(10, 3).Value = 10.09779095434224
$ws.Cells.Item(10, 4).Value = 6.461644001687021
$ws.Cells.Item(10, 5).Value = 11.05840387700759
$ws.Cells.Item(10, 7).Value = 3.640753973742433
$ws.Cells.Item(10, 12).Value = 10.02494121936246
$ws.Cells.Item(10, 14).Value = 17.68567183029592
$ws.Cells.Item(10, 15).Value = 25.79453113809845
$ws.Cells.Item(11, 2).Value = 20.49626152144512
$ws.Cells.Item(11, 3).Value = 10.28570061204379
$ws.Cells.Item(11, 4).Value = 6.558849212125477
$ws.Cells.Item(11, 5).Value = 11.03884441432294
$ws.Cells.Item(11, 7).Value = 3.638913808821585
$ws.Cells.Item(11, 12).Value = 10.05291459391588
$ws.Cells.Item(11, 14).Value = 17.639106121419
$ws.Cells.Item(11, 15).Value = 25.8280642285554
$ws.Cells.Item(12, 2).Value = 20.63746301707219
$ws.Cells.Item(12, 3).Value = 10.35576647966658
$ws.Cells.Item(12, 4).Value = 6.595535472731903
$ws.Cells.Item(12, 5).Value = 11.03170287691598
$ws.Cells.Item(12, 7).Value = 3.638229987688256
$ws.Cells.Item(12, 12).Value = 10.06376123798956
$ws.Cells.Item(12, 14).Value = 17.62175579818257
$ws.Cells.Item(12, 15).Value = 25.84186215594564
$ws.Cells.Item(13, 2).Value = 20.60711422974389
$ws.Cells.Item(13, 3).Value = 10.34072564961583
$ws.Cells.Item(13, 4).Value = 6.58764051977742
$ws.Cells.Item(13, 5).Value = 11.03322913544335
$ws.Cells.Item(13, 7).Value = 3.638376683337103
$ws.Cells.Item(13, 12).Value = 10.06141402946045
$ws.Cells.Item(13, 14).Value = 17.62547993203437
$ws.Cells.Item(13, 15).Value = 25.83884168679857
$ws.Cells.Item(14, 2).Value = 20.50790577215614
$ws.Cells.Item(14, 3).Value = 10.29148705197628
$ws.Cells.Item(14, 4).Value = 6.561870102283357
$ws.Cells.Item(14, 5).Value = 11.0382515590958
$ws.Cells.Item(14, 7).Value = 3.638857290123629
$ws.Cells.Item(14, 12).Value = 10.05380190280932
$ws.Cells.Item(14, 14).Value = 17.63767303415146
$ws.Cells.Item(14, 15).Value = 25.82917737794538
$ws.Cells.Item(15, 2).Value = 20.44695988547787
$ws.Cells.Item(15, 3).Value = 10.26118378503595
$ws.Cells.Item(15, 4).Value = 6.54606778761625
$ws.Cells.Item(15, 5).Value = 11.04136248416509
$ws.Cells.Item(15, 7).Value = 3.639153367758027
$ws.Cells.Item(15, 12).Value = 10.04917212252842
$ws.Cells.Item(15, 14).Value = 17.6451784877092
$ws.Cells.Item(15, 15).Value = 25.8234007935888
$ws.Cells.Item(16, 2).Value = 20.09507622206002
$ws.Cells.Item(16, 3).Value = 10.08535987621001
$ws.Cells.Item(16, 4).Value = 6.455276276250774
$ws.Cells.Item(16, 5).Value = 11.05971923131919
$ws.Cells.Item(16, 7).Value = 3.640876057189958
$ws.Cells.Item(16, 12).Value = 10.02314892113352
$ws.Cells.Item(16, 14).Value = 17.68875470100517
$ws.Cells.Item(16, 15).Value = 25.79249377491926
$ws.Cells.Item(17, 2).Value = 19.8769730493767
$ws.Cells.Item(17, 3).Value = 9.975593272341273
$ws.Cells.Item(17, 4).Value = 6.399400831426524
$ws.Cells.Item(17, 5).Value = 11.07145258471903
$ws.Cells.Item(17, 7).Value = 3.641956118758502
$ws.Cells.Item(17, 12).Value = 10.00764272666764
$ws.Cells.Item(17, 14).Value = 17.71599303611371
$ws.Cells.Item(17, 15).Value = 25.77549611531478
$ws.Cells.Item(18, 2).Value = 19.75075023781397
$ws.Cells.Item(18, 3).Value = 9.911771132809696
$ws.Cells.Item(18, 4).Value = 6.367210254175776
$ws.Cells.Item(18, 5).Value = 11.07837470293052
$ws.Cells.Item(18, 7).Value = 3.642585907851785
$ws.Cells.Item(18, 12).Value = 9.998894189547974
$ws.Cells.Item(18, 14).Value = 17.73184611730078
$ws.Cells.Item(18, 15).Value = 25.76644186942115
$ws.Cells.Item(19, 2).Value = 19.70788458747102
$ws.Cells.Item(19, 3).Value = 9.89004524882828
$ws.Cells.Item(19, 4).Value = 6.356303301518771
$ws.Cells.Item(19, 5).Value = 11.08074818732531
$ws.Cells.Item(19, 7).Value = 3.642800617055949
$ws.Cells.Item(19, 12).Value = 9.995961489466373
$ws.Cells.Item(19, 14).Value = 17.73724574058667
$ws.Cells.Item(19, 15).Value = 25.76350040591263
$ws.Cells.Item(20, 2).Value = 19.90027182645915
$ws.Cells.Item(20, 3).Value = 9.987349538245288
$ws.Cells.Item(20, 4).Value = 6.405354594894716
$ws.Cells.Item(20, 5).Value = 11.07018560065706
$ws.Cells.Item(20, 7).Value = 3.641840258303296
$ws.Cells.Item(20, 12).Value = 10.00927581140683
$ws.Cells.Item(20, 14).Value = 17.71307419656018
$ws.Cells.Item(20, 15).Value = 25.77723080442497
$ws.Cells.Item(21, 2).Value = 20.53708299397731
$ws.Cells.Item(21, 3).Value = 10.30597952786146
$ws.Cells.Item(21, 4).Value = 6.569443148668483
$ws.Cells.Item(21, 5).Value = 11.03676915179049
$ws.Cells.Item(21, 7).Value = 3.638715771753919
$ws.Cells.Item(21, 12).Value = 10.05603093019963
$ws.Cells.Item(21, 14).Value = 17.63408395238733
$ws.Cells.Item(21, 15).Value = 25.83198620926658
$ws.Cells.Item(22, 2).Value = 20.94543554581456
$ws.Cells.Item(22, 3).Value = 10.50784601182544
$ws.Cells.Item(22, 4).Value = 6.67594927993144
$ws.Cells.Item(22, 5).Value = 11.01647551882548
$ws.Cells.Item(22, 7).Value = 3.63674953511708
$ws.Cells.Item(22, 12).Value = 10.08806415358011
$ws.Cells.Item(22, 14).Value = 17.58410874343071
$ws.Cells.Item(22, 15).Value = 25.87417919809993
$ws.Cells.Item(23, 2).Value = 20.72825026572563
$ws.Cells.Item(23, 3).Value = 10.40070112827279
$ws.Cells.Item(23, 4).Value = 6.619184765692003
$ws.Cells.Item(23, 5).Value = 11.02716507134855
$ws.Cells.Item(23, 7).Value = 3.637792040690115
$ws.Cells.Item(23, 12).Value = 10.07083434889517
$ws.Cells.Item(23, 14).Value = 17.61063098667409
$ws.Cells.Item(23, 15).Value = 25.85107524737335
$ws.Cells.Item(24, 2).Value = 19.88974103131321
$ws.Cells.Item(24, 3).Value = 9.982036757292525
$ws.Cells.Item(24, 4).Value = 6.402663104242731
$ws.Cells.Item(24, 5).Value = 11.07075785452475
$ws.Cells.Item(24, 7).Value = 3.64189261125007
$ws.Cells.Item(24, 12).Value = 10.00853697573215
$ws.Cells.Item(24, 14).Value = 17.71439320133883
$ws.Cells.Item(24, 15).Value = 25.77644431481546
$ws.Cells.Item(25, 2).Value = 18.95629179661424
$ws.Cells.Item(25, 3).Value = 9.504356630835371
$ws.Cells.Item(25, 4).Value = 6.167279598646601
$ws.Cells.Item(25, 5).Value = 11.12440824784999
$ws.Cells.Item(25, 7).Value = 3.646643548440361
$ws.Cells.Item(25, 12).Value = 9.947866568177682
$ws.Cells.Item(25, 14).Value = 17.83344675074333
$ws.Cells.Item(25, 15).Value = 25.72356076127379
